# Weekly price-sheet update: a new weekly record for
# "Agrícola del Norte S.A. de Arica - Cebollín baby" (dated 2022-03-21,
# serial 44641) is inserted as the new row 66, pushing the previously
# existing rows 66-78 down to 67-79.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 66 (shifts rows 66:78 -> 67:79, copying the
# formatting of the row that used to be there, same as Excel's own
# "Insert Sheet Rows" command anchored on row 66).
$ws.Rows("66:66").Insert()

# Populate the new row 66 with the new week's record.
$ws.Cells.Item(66, 1).Value2  = 1
$ws.Cells.Item(66, 2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(66, 3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(66, 4).Value2  = 44641
$ws.Cells.Item(66, 5).Value2  = 15
$ws.Cells.Item(66, 6).Value2  = 100112038
$ws.Cells.Item(66, 7).Value2  = "Cebollín baby"
$ws.Cells.Item(66, 8).Value2  = "Sin especificar"
$ws.Cells.Item(66, 9).Value2  = "Primera"
$ws.Cells.Item(66, 10).Value2 = 300
$ws.Cells.Item(66, 11).Value2 = 1300
$ws.Cells.Item(66, 12).Value2 = 1500
$ws.Cells.Item(66, 13).Value2 = 1400
$ws.Cells.Item(66, 14).Value2 = "`$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(66, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(66, 16).Value2 = 700
$ws.Cells.Item(66, 17).Value2 = 2
$ws.Cells.Item(66, 18).Value2 = "Hortaliza"
